# ------------------------------------------------------------------
# Applies the 10/22/2023 4:03PM EST update to the Illegal War Act
# Prevention Security Systems document.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# Track changes would wrap every edit below in <w:ins>/<w:del> markup,
# which is not what the target diff shows - make sure it is off.
$d.TrackRevisions = $false

$wNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-WordXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><?mso-application progid="Word.Document"?><w:wordDocument ' + $wNS + '><w:body>' + $bodyXml + '</w:body></w:wordDocument>'
}

# Replace the text running from the first occurrence of $anchorText inside
# paragraph $paraIndex through to the end of that paragraph's text (but
# not its paragraph mark) with $bodyXml (one or more <w:r>/<w:proofErr>
# elements). Using the "through end of paragraph" range (rather than just
# the exact substring) keeps any untouched trailing runs in their correct
# relative order.
function Replace-ParaTail($doc, [int]$paraIndex, [string]$anchorText, [string]$bodyXml) {
    $p = $doc.Paragraphs.Item($paraIndex)
    $full = $p.Range.Text
    $idx = $full.IndexOf($anchorText)
    if ($idx -lt 0) {
        throw "Replace-ParaTail: anchor [$anchorText] not found in paragraph $paraIndex : [$full]"
    }
    $pStart = $p.Range.Start
    $rStart = $pStart + $idx
    $rEnd = $pStart + $full.Length - 1   # stop before the paragraph mark
    $rng = $doc.Range($rStart, $rEnd)
    $rng.InsertXML((New-WordXml $bodyXml))
}

# 1. Header date/time field -------------------------------------------------
$d.Content.Find.Execute("10/22/2023 1:43:07 PM", $true, $false, $false, $false, $false, $true, 1, $false, "10/22/2023 2:49:54 PM", 2) | Out-Null

# 2. "ANY GUERRILLA OPPOSITION GROUP;" -> "ANY GENOCIDE;" --------------------
#    (split into "G" / "ENOCIDE" runs, wrapped in proofErr gramStart/gramEnd)
Replace-ParaTail $d 32 "GUERRILLA OPPOSITION GROUP" (
  '<w:p>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>G</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>ENOCIDE</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>'
)

# 3. "ANY HISTORICAL CATACLYSMS;" -> "ANY GUERRILLA OPPOSITION GROUP;" ------
#    (split into "GUERRILLA OPPOSITION " / "GROUP" runs, proofErr wrapped)
Replace-ParaTail $d 33 "HISTORICAL CATACLYSMS" (
  '<w:p>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">GUERRILLA OPPOSITION </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>GROUP</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>'
)

# 4. "ANY HORRENDOUS ATROCITY;" -> "ANY HISTORICAL CATACLYSMS;" -------------
Replace-ParaTail $d 34 "HORRENDOUS ATROCITY" (
  '<w:p>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>HISTORICAL CATACLYSMS</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r>' +
  '</w:p>'
)

# 5. "ANY HUMAN RIGHTS ABUSE;" -> "ANY HORRENDOUS ATROCITY;" ----------------
Replace-ParaTail $d 35 "HUMAN RIGHTS ABUSE" (
  '<w:p>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>HORRENDOUS ATROCITY</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r>' +
  '</w:p>'
)

# 6. "ANY HUMANITARIAN ABUSE;" -> "ANY HUMAN RIGHTS ABUSE;" -----------------
Replace-ParaTail $d 36 "HUMANITARIAN ABUSE" (
  '<w:p>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>HUMAN RIGHTS ABUSE</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r>' +
  '</w:p>'
)

# 7. New paragraph inserted: "PREVENTION SECURITY SYSTEM: ANY HUMANITARIAN ABUSE;"
#    goes right after paragraph 36 (the one that now reads ".. HUMAN RIGHTS ABUSE;")
#    and right before the existing ".. ILLEGAL DOMESTIC ACTION;" paragraph.
$p36 = $d.Paragraphs.Item(36)
$p36.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item(37)
$newPara.Range.InsertXML((New-WordXml (
  '<w:p>' +
    '<w:pPr><w:ind w:left="360"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>PREVENTION SECURITY SYSTEM</w:t></w:r>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">ANY </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>HUMANITARIAN ABUSE</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r>' +
  '</w:p>'
)))

# 8. Merge the " " + "          " (1 + 10 space) runs between "... OR" and
#    " ANY LARGE-SCALE ..." into a single 11-space run.
$paraIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "LARGE-SCALE CRIMES AGAINST CIVILIZATION") {
        $paraIdx = $i
        break
    }
}
if ($paraIdx -eq -1) { throw "Could not find the SYSTEMATIC/LARGE-SCALE paragraph" }
$p = $d.Paragraphs.Item($paraIdx)
$full = $p.Range.Text
$orIdx = $full.IndexOf(" OR")
$start = $orIdx + (" OR").Length
$pStart = $p.Range.Start
$rStart = $pStart + $start
$rEnd = $pStart + $full.Length - 1
$rng = $d.Range($rStart, $rEnd)
$rng.InsertXML((New-WordXml (
  '<w:p>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">           </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> ANY </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>LARGE-SCALE CRIMES AGAINST CIVILIZATION</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r>' +
  '</w:p>'
)))

# 9. Add <w:lastRenderedPageBreak/> right before the final "}" run ----------
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$full = $lastP.Range.Text
$pStart = $lastP.Range.Start
$rStart = $pStart
$rEnd = $pStart + $full.Length - 1
$rng = $d.Range($rStart, $rEnd)
$rng.InsertXML((New-WordXml (
  '<w:p><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>}</w:t></w:r></w:p>'
)))
